# Generate Report for Handoff
# This script mutates the localization-status workbook so that the three
# files that were previously "Handed back" (292c40e8... and d7adbe74...)
# are replaced by the new handoff-tracking data: two images that are
# "Ready for handoff" (4552959a...png and 460ffba4...png) plus a markdown
# file (64eaf0e6...md) that is a dependency of those two images and is
# itself included/ready as well.

$wb = $excel.ActiveWorkbook

$commit = "db9842e00bc1f7d63fd6dd4afd8d5bba261a0d04"
$base = "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e"

$img1 = "4552959a-6d1b-4951-82ff-899c298499db.png"
$img2 = "460ffba4-ee0b-48eb-9fe0-cdb05fff5c3e.png"
$md   = "64eaf0e6-3d2e-4d69-bcec-d168f078bc58.md"

$img1Target = "b27a0ce97553d72628356f7560fb3ba156025d3f.png"
$img2Target = "03c3e703f8b1f61806d723d4665548424d26a825.png"
$mdZhTarget = "64eaf0e6-3d2e-4d69-bcec-d168f078bc58.cd453a22425f76a767978edeb2093488a81e957b.zh-cn.xlf"
$mdDeTarget = "64eaf0e6-3d2e-4d69-bcec-d168f078bc58.cd453a22425f76a767978edeb2093488a81e957b.de-de.xlf"

$status = "Ready for handoff"
$dateZh = "2016-03-18 17:22:57"
$dateDe = "2016-03-18 17:23:08"
$dateOverview = "2016-03-18 17:23:08"
$epoch = "0001-01-01 00:00:00"
$depFrom = "e2e\$md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = $img1
$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = $dateOverview

$ws1.Range("A3").Value = $img2
$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = $dateOverview

$ws1.Range("A4").Value = $md
$ws1.Range("B4").Value = $status
$ws1.Range("C4").Value = $status
$ws1.Range("D4").Value = $dateOverview

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$base/$img1", "", "", $img1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$base/$img2", "", "", $img2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$base/$md", "", "", $md)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("F2:G2").ClearContents()
$ws2.Range("F2:G2").ClearFormats()
$ws2.Range("I2").ClearContents()

$ws2.Range("F3:G3").ClearContents()
$ws2.Range("F3:G3").ClearFormats()
$ws2.Range("I3").ClearContents()

$ws2.Range("A2").Value = $img1
$ws2.Range("B2").Value = ".png"
$ws2.Range("C2").Value = $status
$ws2.Range("D2").Value = $img1Target
$ws2.Range("E2").Value = $dateZh
$ws2.Range("H2").Value = $epoch
$ws2.Range("J2").Value = "IsDependency"
$ws2.Range("K2").Value = $depFrom

$ws2.Range("A3").Value = $img2
$ws2.Range("B3").Value = ".png"
$ws2.Range("C3").Value = $status
$ws2.Range("D3").Value = $img2Target
$ws2.Range("E3").Value = $dateZh
$ws2.Range("H3").Value = $epoch
$ws2.Range("J3").Value = "IsDependency"
$ws2.Range("K3").Value = $depFrom

$ws2.Range("A4").Value = $md
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = $status
$ws2.Range("D4").Value = $mdZhTarget
$ws2.Range("E4").Value = $dateZh
$ws2.Range("H4").Value = $epoch
$ws2.Range("I4").Value = ""
$ws2.Range("J4").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$base/$img1", "", "", $img1)
$ws2.Hyperlinks.Add($ws2.Range("B2"), "$base/$img1", "", "", ".png")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad745b0dbfebfeea2549edf6a4c7b5e58a80baf1/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$img1Target", "", "", $img1Target)

$ws2.Hyperlinks.Add($ws2.Range("A3"), "$base/$img2", "", "", $img2)
$ws2.Hyperlinks.Add($ws2.Range("B3"), "$base/$img2", "", "", ".png")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad745b0dbfebfeea2549edf6a4c7b5e58a80baf1/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$img2Target", "", "", $img2Target)

$ws2.Hyperlinks.Add($ws2.Range("A4"), "$base/$md", "", "", $md)
$ws2.Hyperlinks.Add($ws2.Range("B4"), "$base/$md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad745b0dbfebfeea2549edf6a4c7b5e58a80baf1/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$mdZhTarget", "", "", $mdZhTarget)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("F2:G2").ClearContents()
$ws3.Range("F2:G2").ClearFormats()
$ws3.Range("I2").ClearContents()

$ws3.Range("F3:G3").ClearContents()
$ws3.Range("F3:G3").ClearFormats()
$ws3.Range("I3").ClearContents()

$ws3.Range("A2").Value = $img1
$ws3.Range("B2").Value = ".png"
$ws3.Range("C2").Value = $status
$ws3.Range("D2").Value = $img1Target
$ws3.Range("E2").Value = $dateDe
$ws3.Range("H2").Value = $epoch
$ws3.Range("J2").Value = "IsDependency"
$ws3.Range("K2").Value = $depFrom

$ws3.Range("A3").Value = $img2
$ws3.Range("B3").Value = ".png"
$ws3.Range("C3").Value = $status
$ws3.Range("D3").Value = $img2Target
$ws3.Range("E3").Value = $dateDe
$ws3.Range("H3").Value = $epoch
$ws3.Range("J3").Value = "IsDependency"
$ws3.Range("K3").Value = $depFrom

$ws3.Range("A4").Value = $md
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = $status
$ws3.Range("D4").Value = $mdDeTarget
$ws3.Range("E4").Value = $dateDe
$ws3.Range("H4").Value = $epoch
$ws3.Range("I4").Value = ""
$ws3.Range("J4").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$base/$img1", "", "", $img1)
$ws3.Hyperlinks.Add($ws3.Range("B2"), "$base/$img1", "", "", ".png")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73997616a69e58ad46d1b551185887b0c8226387/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$img1Target", "", "", $img1Target)

$ws3.Hyperlinks.Add($ws3.Range("A3"), "$base/$img2", "", "", $img2)
$ws3.Hyperlinks.Add($ws3.Range("B3"), "$base/$img2", "", "", ".png")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73997616a69e58ad46d1b551185887b0c8226387/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$img2Target", "", "", $img2Target)

$ws3.Hyperlinks.Add($ws3.Range("A4"), "$base/$md", "", "", $md)
$ws3.Hyperlinks.Add($ws3.Range("B4"), "$base/$md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73997616a69e58ad46d1b551185887b0c8226387/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$mdDeTarget", "", "", $mdDeTarget)

Write-Host "Report regenerated for handoff."
